$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 4: clear the old status cell (F4 had "DONE")
$ws.Range("F4").ClearContents()

# Row 5: add priority value + new "IN PROGRESS" status
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "IN PROGRESS"

# Row 6: add priority value + new "IN PROGRESS" status
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "IN PROGRESS"

# Update the active selection to mirror the saved view state
$ws.Activate()
$ws.Range("E7").Select()
